$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text in this workbook (note the
# thousands-separator dots rather than true numeric formatting, e.g. "26.053.66").
# Assigning a numeric-looking string via .Value lets Excel auto-detect it as a
# Number, which would silently drop significant trailing zeros (e.g. "4.500" ->
# 4.5) or flip to scientific notation (e.g. "0.000008502" -> 8.502E-06). Force
# those specific cells to Text format first so the literal string is preserved.

$ws.Range("D2").Value = "26.053.66"
$ws.Range("E2").Value = "  -1.94%  "
$ws.Range("D3").Value = "1.665.78"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.43"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5098"
$ws.Range("E6").Value = "  +1.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2624"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06407"
$ws.Range("E9").Value = "  +3.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.65"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07424"
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("D12").Value = "1.668.16"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.500"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5802"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008502"
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.23"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "26.124.31"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.903"
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.74"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.23"
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.194"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.63"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.612"
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1188"
$ws.Range("E26").Value = "  +4.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.58"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06434"
$ws.Range("E28").Value = "  +12.36%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.315"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.502"
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.628"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6051"
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.685"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("E38").Value = "  +5.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01610"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").Value = "1.072.90"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8585"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.59"
$ws.Range("E43").Value = "  +2.22%  "
$ws.Range("D44").Value = "1.813.65"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("E45").Value = "  +6.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.06"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05207"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4294"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.934"
$ws.Range("E51").Value = "  +4.89%  "
